$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 60, shifting the existing rows 60-166 down to 61-167.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new weekly record.
$ws.Range("A60").Value = 10
$ws.Range("B60").Value = "Vega Modelo de Temuco"
$ws.Range("C60").Value = "La Araucanía"
$ws.Range("D60").Value = 44571
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = 100112043
$ws.Range("G60").Value = "Pepino dulce"
$ws.Range("H60").Value = "Cultivar XV región"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 35
$ws.Range("K60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("M60").Value = 25000
$ws.Range("N60").Value = "`$/bandeja 18 kilos"
$ws.Range("O60").Value = "Provincia de Limarí"
$ws.Range("P60").Value = 1389
$ws.Range("Q60").Value = 18
$ws.Range("R60").Value = "Hortaliza"
